$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.20950566666667
$ws.Range("H2").Value = 87.628517
$ws.Range("I2").Value = 0.01829497698069002
$ws.Range("J2").Value = 0.01840828041918582
$ws.Range("M2").Value = 14.349718
$ws.Range("N2").Value = 43.049154
$ws.Range("O2").Value = 0.1016415840981481
$ws.Range("P2").Value = 0.1034081666702025
$ws.Range("Q2").Value = 419.1481692360687
$ws.Range("R2").Value = 3772.333523124618
$ws.Range("S2").Value = 0.001859530441356488
$ws.Range("T2").Value = 0.001903566529698993
$ws.Range("G3").Value = 29.20950566666667
$ws.Range("H3").Value = 87.628517
$ws.Range("I3").Value = 0.01829497698069002
$ws.Range("J3").Value = 0.01840828041918582
$ws.Range("O3").Value = 0.04778708884009916
$ws.Range("P3").Value = 0.04861765281706964
$ws.Range("Q3").Value = 197.063740969518
$ws.Range("R3").Value = 1773.573668725662
$ws.Range("S3").Value = 0.000874263690303803
$ws.Range("T3").Value = 0.0008949673863792374
$ws.Range("G4").Value = 29.20950566666667
$ws.Range("H4").Value = 87.628517
$ws.Range("I4").Value = 0.01829497698069002
$ws.Range("J4").Value = 0.01840828041918582
$ws.Range("M4").Value = 66.43651233333334
$ws.Range("N4").Value = 199.309537
$ws.Range("O4").Value = 0.4705815372480596
$ws.Range("P4").Value = 0.4787604843769264
$ws.Range("Q4").Value = 1940.57768347407
$ws.Range("R4").Value = 17465.19915126663
$ws.Range("S4").Value = 0.008609278391490973
$ws.Range("T4").Value = 0.008813157250035692
$ws.Range("G5").Value = 29.20950566666667
$ws.Range("H5").Value = 87.628517
$ws.Range("I5").Value = 0.01829497698069002
$ws.Range("J5").Value = 0.01840828041918582
$ws.Range("M5").Value = 7.2355625
$ws.Range("N5").Value = 14.471125
$ws.Range("O5").Value = 0.05125076564857627
$ws.Range("P5").Value = 0.03476102006337534
$ws.Range("Q5").Value = 211.3472038452709
$ws.Range("R5").Value = 1268.083223071625
$ws.Range("S5").Value = 0.0009376315777834416
$ws.Range("T5").Value = 0.0006398906049835577
$ws.Range("G6").Value = 29.20950566666667
$ws.Range("H6").Value = 87.628517
$ws.Range("I6").Value = 0.01829497698069002
$ws.Range("J6").Value = 0.01840828041918582
$ws.Range("M6").Value = 46.41124333333334
$ws.Range("N6").Value = 139.23373
$ws.Range("O6").Value = 0.328739024165117
$ws.Range("P6").Value = 0.3344526760724259
$ws.Range("Q6").Value = 1355.649475142046
$ws.Range("R6").Value = 12200.84527627841
$ws.Range("S6").Value = 0.006014272879755315
$ws.Range("T6").Value = 0.006156698648088336
$ws.Range("I7").Value = 0.913374480506715
$ws.Range("J7").Value = 0.9190311407684336
$ws.Range("M7").Value = 14.349718
$ws.Range("N7").Value = 43.049154
$ws.Range("O7").Value = 0.1016415840981481
$ws.Range("P7").Value = 0.1034081666702025
$ws.Range("Q7").Value = 20925.92090907871
$ws.Range("R7").Value = 188333.2881817084
$ws.Range("S7").Value = 0.09283682907352557
$ws.Range("T7").Value = 0.09503532537968855
$ws.Range("I8").Value = 0.913374480506715
$ws.Range("J8").Value = 0.9190311407684336
$ws.Range("O8").Value = 0.04778708884009916
$ws.Range("P8").Value = 0.04861765281706964
$ws.Range("S8").Value = 0.0436475074442538
$ws.Range("T8").Value = 0.04468113692995516
$ws.Range("I9").Value = 0.913374480506715
$ws.Range("J9").Value = 0.9190311407684336
$ws.Range("M9").Value = 66.43651233333334
$ws.Range("N9").Value = 199.309537
$ws.Range("O9").Value = 0.4705815372480596
$ws.Range("P9").Value = 0.4787604843769264
$ws.Range("Q9").Value = 96883.10268970899
$ws.Range("R9").Value = 871947.9242073809
$ws.Range("S9").Value = 0.4298171671199978
$ws.Range("T9").Value = 0.4399957941117745
$ws.Range("I10").Value = 0.913374480506715
$ws.Range("J10").Value = 0.9190311407684336
$ws.Range("M10").Value = 7.2355625
$ws.Range("N10").Value = 14.471125
$ws.Range("O10").Value = 0.05125076564857627
$ws.Range("P10").Value = 0.03476102006337534
$ws.Range("Q10").Value = 10551.48321435277
$ws.Range("R10").Value = 63308.89928611663
$ws.Range("S10").Value = 0.04681114144983974
$ws.Range("T10").Value = 0.03194645992311825
$ws.Range("I11").Value = 0.913374480506715
$ws.Range("J11").Value = 0.9190311407684336
$ws.Range("M11").Value = 46.41124333333334
$ws.Range("N11").Value = 139.23373
$ws.Range("O11").Value = 0.328739024165117
$ws.Range("P11").Value = 0.3344526760724259
$ws.Range("Q11").Value = 67680.63367414885
$ws.Range("R11").Value = 609125.7030673396
$ws.Range("S11").Value = 0.3002618354190981
$ws.Range("T11").Value = 0.307372424423897
$ws.Range("G12").Value = 57.98602933333333
$ws.Range("H12").Value = 173.958088
$ws.Range("I12").Value = 0.03631876156896331
$ws.Range("J12").Value = 0.03654368891224535
$ws.Range("M12").Value = 14.349718
$ws.Range("N12").Value = 43.049154
$ws.Range("O12").Value = 0.1016415840981481
$ws.Range("P12").Value = 0.1034081666702025
$ws.Range("Q12").Value = 832.0831688730614
$ws.Range("R12").Value = 7488.748519857551
$ws.Range("S12").Value = 0.003691496458352372
$ws.Range("T12").Value = 0.003778915873781499
$ws.Range("G13").Value = 57.98602933333333
$ws.Range("H13").Value = 173.958088
$ws.Range("I13").Value = 0.03631876156896331
$ws.Range("J13").Value = 0.03654368891224535
$ws.Range("O13").Value = 0.04778708884009916
$ws.Range("P13").Value = 0.04861765281706964
$ws.Range("Q13").Value = 391.206342031152
$ws.Range("R13").Value = 3520.857078280367
$ws.Range("S13").Value = 0.001735567885658429
$ws.Range("T13").Value = 0.001776668380190542
$ws.Range("G14").Value = 57.98602933333333
$ws.Range("H14").Value = 173.958088
$ws.Range("I14").Value = 0.03631876156896331
$ws.Range("J14").Value = 0.03654368891224535
$ws.Range("M14").Value = 66.43651233333334
$ws.Range("N14").Value = 199.309537
$ws.Range("O14").Value = 0.4705815372480596
$ws.Range("P14").Value = 0.4787604843769264
$ws.Range("Q14").Value = 3852.389552965029
$ws.Range("R14").Value = 34671.50597668525
$ws.Range("S14").Value = 0.01709093865006851
$ws.Range("T14").Value = 0.0174956742045463
$ws.Range("G15").Value = 57.98602933333333
$ws.Range("H15").Value = 173.958088
$ws.Range("I15").Value = 0.03631876156896331
$ws.Range("J15").Value = 0.03654368891224535
$ws.Range("M15").Value = 7.2355625
$ws.Range("N15").Value = 14.471125
$ws.Range("O15").Value = 0.05125076564857627
$ws.Range("P15").Value = 0.03476102006337534
$ws.Range("Q15").Value = 419.5615393681666
$ws.Range("R15").Value = 2517.369236209
$ws.Range("S15").Value = 0.001861364337817457
$ws.Range("T15").Value = 0.001270295903468307
$ws.Range("G16").Value = 57.98602933333333
$ws.Range("H16").Value = 173.958088
$ws.Range("I16").Value = 0.03631876156896331
$ws.Range("J16").Value = 0.03654368891224535
$ws.Range("M16").Value = 46.41124333333334
$ws.Range("N16").Value = 139.23373
$ws.Range("O16").Value = 0.328739024165117
$ws.Range("P16").Value = 0.3344526760724259
$ws.Range("Q16").Value = 2691.203717323138
$ws.Range("R16").Value = 24220.83345590824
$ws.Range("S16").Value = 0.01193939423706655
$ws.Range("T16").Value = 0.0122221345502587
$ws.Range("G17").Value = 29.481085
$ws.Range("H17").Value = 58.96217
$ws.Range("I17").Value = 0.01846507700595112
$ws.Range("J17").Value = 0.01238628926567028
$ws.Range("M17").Value = 14.349718
$ws.Range("N17").Value = 43.049154
$ws.Range("O17").Value = 0.1016415840981481
$ws.Range("P17").Value = 0.1034081666702025
$ws.Range("Q17").Value = 423.04525608403
$ws.Range("R17").Value = 2538.27153650418
$ws.Range("S17").Value = 0.001876819677379161
$ws.Range("T17").Value = 0.001280843464809773
$ws.Range("G18").Value = 29.481085
$ws.Range("H18").Value = 58.96217
$ws.Range("I18").Value = 0.01846507700595112
$ws.Range("J18").Value = 0.01238628926567028
$ws.Range("O18").Value = 0.04778708884009916
$ws.Range("P18").Value = 0.04861765281706964
$ws.Range("Q18").Value = 198.89596777977
$ws.Range("R18").Value = 1193.37580667862
$ws.Range("S18").Value = 0.0008823922753226585
$ws.Range("T18").Value = 0.000602192311210154
$ws.Range("G19").Value = 29.481085
$ws.Range("H19").Value = 58.96217
$ws.Range("I19").Value = 0.01846507700595112
$ws.Range("J19").Value = 0.01238628926567028
$ws.Range("M19").Value = 66.43651233333334
$ws.Range("N19").Value = 199.309537
$ws.Range("O19").Value = 0.4705815372480596
$ws.Range("P19").Value = 0.4787604843769264
$ws.Range("Q19").Value = 1958.620467202549
$ws.Range("R19").Value = 11751.72280321529
$ws.Range("S19").Value = 0.008689324322864279
$ws.Range("T19").Value = 0.005930065848465026
$ws.Range("G20").Value = 29.481085
$ws.Range("H20").Value = 58.96217
$ws.Range("I20").Value = 0.01846507700595112
$ws.Range("J20").Value = 0.01238628926567028
$ws.Range("M20").Value = 7.2355625
$ws.Range("N20").Value = 14.471125
$ws.Range("O20").Value = 0.05125076564857627
$ws.Range("P20").Value = 0.03476102006337534
$ws.Range("Q20").Value = 213.3122330853125
$ws.Range("R20").Value = 853.24893234125
$ws.Range("S20").Value = 0.0009463493343149154
$ws.Range("T20").Value = 0.0004305600496747351
$ws.Range("G21").Value = 29.481085
$ws.Range("H21").Value = 58.96217
$ws.Range("I21").Value = 0.01846507700595112
$ws.Range("J21").Value = 0.01238628926567028
$ws.Range("M21").Value = 46.41124333333334
$ws.Range("N21").Value = 139.23373
$ws.Range("O21").Value = 0.328739024165117
$ws.Range("P21").Value = 0.3344526760724259
$ws.Range("Q21").Value = 1368.253809665684
$ws.Range("R21").Value = 8209.522857994101
$ws.Range("S21").Value = 0.006070191396070112
$ws.Range("T21").Value = 0.004142627591510587
$ws.Range("G22").Value = 21.628479
$ws.Range("H22").Value = 64.885437
$ws.Range("I22").Value = 0.01354670393768061
$ws.Range("J22").Value = 0.01363060063446486
$ws.Range("M22").Value = 14.349718
$ws.Range("N22").Value = 43.049154
$ws.Range("O22").Value = 0.1016415840981481
$ws.Range("P22").Value = 0.1034081666702025
$ws.Range("Q22").Value = 310.362574418922
$ws.Range("R22").Value = 2793.263169770298
$ws.Range("S22").Value = 0.001376908447534477
$ws.Range("T22").Value = 0.001409515422223711
$ws.Range("G23").Value = 21.628479
$ws.Range("H23").Value = 64.885437
$ws.Range("I23").Value = 0.01354670393768061
$ws.Range("J23").Value = 0.01363060063446486
$ws.Range("O23").Value = 0.04778708884009916
$ws.Range("P23").Value = 0.04861765281706964
$ws.Range("Q23").Value = 145.917874539198
$ws.Range("R23").Value = 1313.260870852782
$ws.Range("S23").Value = 0.0006473575445604644
$ws.Range("T23").Value = 0.0006626878093345418
$ws.Range("G24").Value = 21.628479
$ws.Range("H24").Value = 64.885437
$ws.Range("I24").Value = 0.01354670393768061
$ws.Range("J24").Value = 0.01363060063446486
$ws.Range("M24").Value = 66.43651233333334
$ws.Range("N24").Value = 199.309537
$ws.Range("O24").Value = 0.4705815372480596
$ws.Range("P24").Value = 0.4787604843769264
$ws.Range("Q24").Value = 1436.920711834741
$ws.Range("R24").Value = 12932.28640651267
$ws.Range("S24").Value = 0.006374828763638084
$ws.Range("T24").Value = 0.006525792962104838
$ws.Range("G25").Value = 21.628479
$ws.Range("H25").Value = 64.885437
$ws.Range("I25").Value = 0.01354670393768061
$ws.Range("J25").Value = 0.01363060063446486
$ws.Range("M25").Value = 7.2355625
$ws.Range("N25").Value = 14.471125
$ws.Range("O25").Value = 0.05125076564857627
$ws.Range("P25").Value = 0.03476102006337534
$ws.Range("Q25").Value = 156.4942115844375
$ws.Range("R25").Value = 938.965269506625
$ws.Range("S25").Value = 0.0006942789488207143
$ws.Range("T25").Value = 0.0004738135821304897
$ws.Range("G26").Value = 21.628479
$ws.Range("H26").Value = 64.885437
$ws.Range("I26").Value = 0.01354670393768061
$ws.Range("J26").Value = 0.01363060063446486
$ws.Range("M26").Value = 46.41124333333334
$ws.Range("N26").Value = 139.23373
$ws.Range("O26").Value = 0.328739024165117
$ws.Range("P26").Value = 0.3344526760724259
$ws.Range("Q26").Value = 1003.80460179889
$ws.Range("R26").Value = 9034.241416190011
$ws.Range("S26").Value = 0.004453330233126871
$ws.Range("T26").Value = 0.00455879085867128
